$wb = $excel.ActiveWorkbook

# --- Rename sheet "getByEmail" -> "products" ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "products"

# --- Build the JSON strings for the three products ---
$mouseJson    = "{`n ""name"":""mouse"",`n ""description"":""mouse_description"",`n ""price"":10.14`n}"
$keyboardJson = "{`n ""name"":""keyBoard"",`n ""description"":""104 button keyboard"",`n ""price"":20.33`n}"
$laptopJson   = "{`n ""name"":""laptop"",`n ""description"":""laptop_description"",`n ""price"":999.99`n}"

# --- Populate the data (row 1 already has a string cell pointing into sharedStrings,
#     overwrite it and add rows 2 and 3) ---
$ws3.Range("A1").Value = $mouseJson
$ws3.Range("A2").Value = $keyboardJson
$ws3.Range("A3").Value = $laptopJson

# --- Apply wrap text + column width + row heights ---
$rng = $ws3.Range("A1:A3")
$rng.WrapText = $true

$ws3.Columns.Item(1).ColumnWidth = 33

$ws3.Rows.Item(1).RowHeight = 75
$ws3.Rows.Item(2).RowHeight = 90
$ws3.Rows.Item(3).RowHeight = 75

# --- Update the selection to L3 without leaving sheet "products" as the active tab ---
$ws1 = $wb.Worksheets.Item(1)
$ws3.Activate() | Out-Null
$ws3.Range("L3").Select() | Out-Null
$ws1.Activate() | Out-Null

Write-Host "done"
